$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shared-string text fixes (row 1 header labels): the Kyrgyz and Russian
#    labels had been swapped relative to their language columns, and all
#    three labels used "16.5.1.1a." (dot) instead of "16.5.1.1a " (space).
#    Re-assigning the values here lets Excel rebuild the shared-string table
#    so A1 ends up pointing at the (corrected) Kyrgyz text and B1 at the
#    (corrected) Russian text, matching the original column headers below
#    them (A = Kyrgyz, B = Russian, C = English).
# ---------------------------------------------------------------------------
$kyrgyzText  = '16.5.1.1a "Аткаруу бийлигинин мамлекеттик органдарындагы жана жергиликтүү өз алдынча башкаруу органдарындагы коррупциянын деңгээли жөнүндө жеке түшүнүк" индекси'
$russianText = "16.5.1.1a Индекс ""Личное представление об уровне коррупции в государственных органах исполнительной власти и органах местного самоуправления''"
$englishText = "16.5.1.1a Index ""Personal views about the level of corruption in executive government authorities and local government''"

$ws.Cells.Item(1,1).Value = $kyrgyzText
$ws.Cells.Item(1,2).Value = $russianText
$ws.Cells.Item(1,3).Value = $englishText

# ---------------------------------------------------------------------------
# 2. Add a new "2020" data column (column I) alongside the existing
#    2015-2019 columns (D-H). Each new cell is formatted like the cell
#    immediately to its left (same row) with a one-decimal "0.0" number
#    format, then given its 2020 value.
# ---------------------------------------------------------------------------
function Set-YearColumnCell($row, $value) {
    $ws.Cells.Item($row, 8).Copy() | Out-Null
    $ws.Cells.Item($row, 9).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 9).Value = $value
}

# Header year
Set-YearColumnCell 4 2020

# Data rows (2020 values)
Set-YearColumnCell 5 12.3
Set-YearColumnCell 6 40.3
Set-YearColumnCell 7 36.2
Set-YearColumnCell 8 44.3
Set-YearColumnCell 9 36
Set-YearColumnCell 10 2.7
Set-YearColumnCell 11 32.9
Set-YearColumnCell 12 11.3
Set-YearColumnCell 13 -18.2
Set-YearColumnCell 14 33

# The new numeric cells (rows 5-14) need the "0.0" number format, unlike the
# plain-integer header cell in row 4.
$ws.Range($ws.Cells.Item(5,9), $ws.Cells.Item(14,9)).NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 3. Move the active selection to reflect the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("F16").Select() | Out-Null
